$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells F1:H1 ---------------------------------------------
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Match the formatting already used by the other header cells (A1:E1):
# bold font, thin box border, centered horizontally, top vertically.
$headerRange = $ws.Range("F1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- New boolean data columns F2:H8 -------------------------------------
$ws.Range("F2").Value = $false
$ws.Range("G2").Value = $false
$ws.Range("H2").Value = $false

$ws.Range("F3").Value = $false
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = $false

$ws.Range("F4").Value = $false
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = $true

$ws.Range("F5").Value = $false
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = $false

$ws.Range("F6").Value = $false
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = $false

$ws.Range("F7").Value = $true
$ws.Range("G7").Value = $true
$ws.Range("H7").Value = $true

$ws.Range("F8").Value = $false
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = $false
